$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(5,3) "149"
Set-TextValue $ws.Cells.Item(5,4) "398886.40"
Set-TextValue $ws.Cells.Item(6,3) "430"
Set-TextValue $ws.Cells.Item(6,4) "1117510.82"
Set-TextValue $ws.Cells.Item(7,3) "155"
Set-TextValue $ws.Cells.Item(7,4) "351041.00"
Set-TextValue $ws.Cells.Item(8,3) "803"
Set-TextValue $ws.Cells.Item(8,4) "3040058.81"
Set-TextValue $ws.Cells.Item(13,3) "97"
Set-TextValue $ws.Cells.Item(13,4) "236800.00"
Set-TextValue $ws.Cells.Item(16,3) "135"
Set-TextValue $ws.Cells.Item(16,4) "606217.26"
Set-TextValue $ws.Cells.Item(37,3) "378"
Set-TextValue $ws.Cells.Item(37,4) "1509848.18"
Set-TextValue $ws.Cells.Item(43,3) "29"
Set-TextValue $ws.Cells.Item(43,4) "78971.00"
Set-TextValue $ws.Cells.Item(52,3) "271"
Set-TextValue $ws.Cells.Item(52,4) "1006390.80"
Set-TextValue $ws.Cells.Item(56,3) "39"
Set-TextValue $ws.Cells.Item(56,4) "100000.00"
Set-TextValue $ws.Cells.Item(60,3) "43"
Set-TextValue $ws.Cells.Item(60,4) "184656.00"
Set-TextValue $ws.Cells.Item(78,3) "212"
Set-TextValue $ws.Cells.Item(78,4) "589693.00"
Set-TextValue $ws.Cells.Item(80,3) "492"
Set-TextValue $ws.Cells.Item(80,4) "2149734.03"
Set-TextValue $ws.Cells.Item(105,3) "14"
Set-TextValue $ws.Cells.Item(105,4) "35909.00"
Set-TextValue $ws.Cells.Item(107,3) "73"
Set-TextValue $ws.Cells.Item(107,4) "180310.00"
Set-TextValue $ws.Cells.Item(108,3) "40"
Set-TextValue $ws.Cells.Item(108,4) "125434.00"
Set-TextValue $ws.Cells.Item(109,3) "15"
Set-TextValue $ws.Cells.Item(109,4) "55913.61"
Set-TextValue $ws.Cells.Item(110,3) "86"
Set-TextValue $ws.Cells.Item(110,4) "532606.82"
Set-TextValue $ws.Cells.Item(111,3) "8"
Set-TextValue $ws.Cells.Item(111,4) "25500.00"
Set-TextValue $ws.Cells.Item(112,3) "5"
Set-TextValue $ws.Cells.Item(112,4) "11500.00"
Set-TextValue $ws.Cells.Item(113,3) "26"
Set-TextValue $ws.Cells.Item(113,4) "83620.00"
Set-TextValue $ws.Cells.Item(114,3) "26"
Set-TextValue $ws.Cells.Item(114,4) "71895.00"
Set-TextValue $ws.Cells.Item(115,3) "14"
Set-TextValue $ws.Cells.Item(115,4) "33100.00"
Set-TextValue $ws.Cells.Item(116,3) "7"
Set-TextValue $ws.Cells.Item(116,4) "17000.00"
Set-TextValue $ws.Cells.Item(118,3) "26"
Set-TextValue $ws.Cells.Item(118,4) "83197.00"
Set-TextValue $ws.Cells.Item(124,3) "498"
Set-TextValue $ws.Cells.Item(124,4) "2240983.06"
Set-TextValue $ws.Cells.Item(133,3) "128"
Set-TextValue $ws.Cells.Item(133,4) "335995.68"
Set-TextValue $ws.Cells.Item(201,3) "664"
Set-TextValue $ws.Cells.Item(201,4) "2541782.58"
Set-TextValue $ws.Cells.Item(202,3) "26"
Set-TextValue $ws.Cells.Item(202,4) "98238.00"
Set-TextValue $ws.Cells.Item(254,3) "130"
Set-TextValue $ws.Cells.Item(254,4) "471264.66"
